$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('O3').Value = 'أهم المميزات إ...'
$ws.Range('O4').Value = 'أهم المميزات I...'
$ws.Range('O8').Value = 'أهم المميزات ل...'
$ws.Range('O10').Value = 'أهم المميزات م...'
$ws.Range('O11').Value = 'أهم المميزات ت...'
$ws.Range('O12').Value = 'أهم المميزات ج...'
$ws.Range('O14').Value = 'أهم المميزات  ...'
$ws.Range('O16').Value = 'أهم المميزات E...'
$ws.Range('O26').Value = 'أهم المميزات  ...'
$ws.Range('O33').Value = 'أهم المميزات F...'
$ws.Range('O34').Value = 'أهم المميزات W...'
$ws.Range('O64').Value = 'أهم المميزات E...'
$ws.Range('O69').Value = 'أهم المميزات T...'
$ws.Range('O80').Value = 'أهم المميزات M...'
$ws.Range('O89').Value = 'أهم المميزات M...'
$ws.Range('P3').Value = 'highlights Exp...'
$ws.Range('P4').Value = 'highlights It ...'
$ws.Range('P8').Value = 'highlights For...'
$ws.Range('P10').Value = 'highlights To ...'
$ws.Range('P11').Value = 'highlights Lux...'
$ws.Range('P12').Value = 'highlights Col...'
$ws.Range('P14').Value = 'highlights Enh...'
$ws.Range('P16').Value = 'highlights Eas...'
$ws.Range('P26').Value = 'highlights   P...'
$ws.Range('P33').Value = 'highlights Fro...'
$ws.Range('P34').Value = 'highlights Wit...'
$ws.Range('P64').Value = 'highlights Eac...'
$ws.Range('P69').Value = 'highlights Thi...'
$ws.Range('P80').Value = 'highlights Mad...'
$ws.Range('P89').Value = 'highlights Mai...'
$ws.Range('Q4').Value = 'وصف It cleans ...'
$ws.Range('Q8').Value = 'وصف أن تركيبة ...'
$ws.Range('Q12').Value = 'وصف يحتوي مقيا...'
$ws.Range('Q16').Value = 'وصف Made from ...'
$ws.Range('Q34').Value = 'وصف The perfec...'
$ws.Range('Q64').Value = 'وصف Indulge in...'
$ws.Range('Q69').Value = 'وصف Ziadah Coa...'
$ws.Range('Q79').Value = 'وصف الكافيين ،...'
$ws.Range('Q80').Value = 'وصف ماء مكربن ...'
$ws.Range('Q87').Value = 'وصف مياه معدني...'
$ws.Range('Q88').Value = 'وصف S Refill P...'
$ws.Range('R4').Value = 'DESCRIPTION It...'
$ws.Range('R8').Value = 'DESCRIPTION De...'
$ws.Range('R12').Value = 'DESCRIPTION Ju...'
$ws.Range('R16').Value = 'DESCRIPTION Ma...'
$ws.Range('R34').Value = 'DESCRIPTION Th...'
$ws.Range('R64').Value = 'DESCRIPTION In...'
$ws.Range('R69').Value = 'DESCRIPTION Zi...'
$ws.Range('R88').Value = 'DESCRIPTION S ...'
$ws.Range('S4').Value = 'حقائق غذائية ا...'
$ws.Range('S16').Value = 'حقائق غذائية ا...'
$ws.Range('S23').Value = 'وصف حقائق غذائ...'
$ws.Range('S25').Value = 'حقائق غذائية ا...'
$ws.Range('S26').Value = 'وصف حقائق غذائ...'
$ws.Range('S30').Value = 'حقائق غذائية ا...'
$ws.Range('S32').Value = 'وصف حقائق غذائ...'
$ws.Range('S34').Value = 'حقائق غذائية ا...'
$ws.Range('S43').Value = 'وصف حقائق غذائ...'
$ws.Range('S46').Value = 'وصف حقائق غذائ...'
$ws.Range('S50').Value = 'وصف حقائق غذائ...'
$ws.Range('S55').Value = 'وصف حقائق غذائ...'
$ws.Range('S64').Value = 'حقائق غذائية ا...'
$ws.Range('S69').Value = 'حقائق غذائية ا...'
$ws.Range('S71').Value = 'وصف حقائق غذائ...'
$ws.Range('S72').Value = 'وصف حقائق غذائ...'
$ws.Range('S74').Value = 'وصف حقائق غذائ...'
$ws.Range('S77').Value = 'وصف حقائق غذائ...'
$ws.Range('S79').Value = 'حقائق غذائية ا...'
$ws.Range('S80').Value = 'حقائق غذائية ا...'
$ws.Range('S87').Value = 'حقائق غذائية ا...'
$ws.Range('S89').Value = 'حقائق غذائية ا...'
$ws.Range('S90').Value = 'وصف حقائق غذائ...'
$ws.Range('T23').Value = 'DESCRIPTION Nu...'
$ws.Range('T25').Value = 'DESCRIPTION Nu...'
$ws.Range('T26').Value = 'DESCRIPTION Nu...'
$ws.Range('T30').Value = 'DESCRIPTION Nu...'
$ws.Range('T32').Value = 'DESCRIPTION Nu...'
$ws.Range('T43').Value = 'DESCRIPTION Nu...'
$ws.Range('T46').Value = 'DESCRIPTION Nu...'
$ws.Range('T50').Value = 'DESCRIPTION Nu...'
$ws.Range('T55').Value = 'DESCRIPTION Nu...'
$ws.Range('T71').Value = 'DESCRIPTION Nu...'
$ws.Range('T72').Value = 'DESCRIPTION Nu...'
$ws.Range('T74').Value = 'DESCRIPTION Nu...'
$ws.Range('T77').Value = 'DESCRIPTION Nu...'
$ws.Range('T79').Value = 'DESCRIPTION Nu...'
$ws.Range('T80').Value = 'DESCRIPTION Nu...'
$ws.Range('T87').Value = 'DESCRIPTION Nu...'
$ws.Range('T90').Value = 'DESCRIPTION Nu...'
$ws.Range('U3').Value = 'معلومات  المكو...'
$ws.Range('U4').Value = 'معلومات  شروط ...'
$ws.Range('U7').Value = 'معلومات  المكو...'
$ws.Range('U8').Value = 'معلومات  المكو...'
$ws.Range('U10').Value = 'معلومات  المكو...'
$ws.Range('U11').Value = 'معلومات  المكو...'
$ws.Range('U14').Value = 'معلومات  المكو...'
$ws.Range('U16').Value = 'معلومات  شروط ...'
$ws.Range('U18').Value = 'معلومات  معلوم...'
$ws.Range('U23').Value = 'معلومات  معلوم...'
$ws.Range('U25').Value = 'معلومات  معلوم...'
$ws.Range('U26').Value = 'معلومات  المكو...'
$ws.Range('U30').Value = 'معلومات  معلوم...'
$ws.Range('U32').Value = 'معلومات  المكو...'
$ws.Range('U33').Value = 'معلومات  المكو...'
$ws.Range('U34').Value = 'معلومات  المكو...'
$ws.Range('U35').Value = 'معلومات  المكو...'
$ws.Range('U43').Value = 'معلومات  المكو...'
$ws.Range('U44').Value = 'معلومات  المكو...'
$ws.Range('U46').Value = 'معلومات  المكو...'
$ws.Range('U63').Value = 'معلومات  المكو...'
$ws.Range('U64').Value = 'معلومات  المكو...'
$ws.Range('U67').Value = 'معلومات  المكو...'
$ws.Range('U69').Value = 'معلومات  المكو...'
$ws.Range('U70').Value = 'معلومات  المكو...'
$ws.Range('U71').Value = 'معلومات  المكو...'
$ws.Range('U73').Value = 'معلومات  المكو...'
$ws.Range('U74').Value = 'معلومات  المكو...'
$ws.Range('U75').Value = 'معلومات  المكو...'
$ws.Range('U77').Value = 'معلومات  المكو...'
$ws.Range('U79').Value = 'معلومات  المكو...'
$ws.Range('U80').Value = 'معلومات  المكو...'
$ws.Range('U82').Value = 'معلومات  المكو...'
$ws.Range('U85').Value = 'معلومات  المكو...'
$ws.Range('U87').Value = 'معلومات  المكو...'
$ws.Range('U89').Value = 'معلومات  المكو...'
$ws.Range('V3').Value = 'INFORMATION In...'
$ws.Range('V4').Value = 'INFORMATION St...'
$ws.Range('V7').Value = 'INFORMATION In...'
$ws.Range('V8').Value = 'INFORMATION In...'
$ws.Range('V10').Value = 'INFORMATION In...'
$ws.Range('V11').Value = 'INFORMATION In...'
$ws.Range('V14').Value = 'INFORMATION In...'
$ws.Range('V16').Value = 'INFORMATION St...'
$ws.Range('V18').Value = 'INFORMATION Al...'
$ws.Range('V23').Value = 'INFORMATION Al...'
$ws.Range('V25').Value = 'INFORMATION Al...'
$ws.Range('V26').Value = 'INFORMATION In...'
$ws.Range('V30').Value = 'INFORMATION Al...'
$ws.Range('V32').Value = 'INFORMATION In...'
$ws.Range('V33').Value = 'INFORMATION In...'
$ws.Range('V34').Value = 'INFORMATION In...'
$ws.Range('V35').Value = 'INFORMATION In...'
$ws.Range('V43').Value = 'INFORMATION In...'
$ws.Range('V44').Value = 'INFORMATION In...'
$ws.Range('V46').Value = 'INFORMATION In...'
$ws.Range('V63').Value = 'INFORMATION In...'
$ws.Range('V64').Value = 'INFORMATION In...'
$ws.Range('V67').Value = 'INFORMATION In...'
$ws.Range('V69').Value = 'INFORMATION In...'
$ws.Range('V70').Value = 'INFORMATION In...'
$ws.Range('V71').Value = 'INFORMATION In...'
$ws.Range('V73').Value = 'INFORMATION In...'
$ws.Range('V74').Value = 'INFORMATION In...'
$ws.Range('V75').Value = 'INFORMATION In...'
$ws.Range('V77').Value = 'INFORMATION In...'
$ws.Range('V79').Value = 'INFORMATION In...'
$ws.Range('V80').Value = 'INFORMATION In...'
$ws.Range('V82').Value = 'INFORMATION In...'
$ws.Range('V85').Value = 'INFORMATION In...'
$ws.Range('V87').Value = 'INFORMATION In...'
$ws.Range('V89').Value = 'INFORMATION In...'
